$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3276.75
$ws.Range("I38").Value = 1036
$ws.Range("K38").Value = 3108
$ws.Range("M38").Value = -2736
$ws.Range("H76").Value = 4318.087
$ws.Range("I76").Value = 3812
$ws.Range("K76").Value = 3812
$ws.Range("M76").Value = -3497
$ws.Range("H79").Value = 4318.087
$ws.Range("I79").Value = 3812
$ws.Range("K79").Value = 3812
$ws.Range("M79").Value = -2720
$ws.Range("H133").Value = 114450
$ws.Range("J133").Value = 114450
$ws.Range("L133").Value = 114450
$ws.Range("N133").Value = -124570
$ws.Range("H137").Value = 11113231
$ws.Range("I137").Value = 2384.75
$ws.Range("K137").Value = 7154.25
$ws.Range("M137").Value = -4604.25
$ws.Range("H138").Value = 11653.68
$ws.Range("I138").Value = 14587.5
$ws.Range("J138").Value = 11531.4375
$ws.Range("K138").Value = 43762.5
$ws.Range("L138").Value = 34594.3125
$ws.Range("M138").Value = -38622.5
$ws.Range("N138").Value = -44874.3125
$ws.Range("H141").Value = 956.6486
$ws.Range("I141").Value = 956.8333
$ws.Range("K141").Value = 2870.4999
$ws.Range("M141").Value = 2309.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4856775.5
$ws.Range("I32").Value = 6367672
$ws.Range("K32").Value = 6367672
$ws.Range("M32").Value = -6367385

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2492.077
$ws.Range("I20").Value = 2260
$ws.Range("J20").Value = 2637.125
$ws.Range("K20").Value = 2260
$ws.Range("L20").Value = 2637.125
$ws.Range("M20").Value = -2013
$ws.Range("N20").Value = -3131.125
$ws.Range("H22").Value = 1189.5
$ws.Range("I22").Value = 1189.5
$ws.Range("K22").Value = 1189.5
$ws.Range("M22").Value = -1016.5
$ws.Range("H94").Value = 3785.818
$ws.Range("I94").Value = 3816.111
$ws.Range("K94").Value = 3816.111
$ws.Range("M94").Value = -3365.111
$ws.Range("H99").Value = 8621.280000000001
$ws.Range("I99").Value = 9365.091
$ws.Range("K99").Value = 9365.091
$ws.Range("M99").Value = -7867.091
$ws.Range("H105").Value = 12828865
$ws.Range("I105").Value = 20842122
$ws.Range("K105").Value = 20842122
$ws.Range("M105").Value = -20840375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2876.4119
$ws.Range("I22").Value = 556.1875
$ws.Range("K22").Value = 556.1875
$ws.Range("M22").Value = -206.1875
$ws.Range("H31").Value = 4629.7646
$ws.Range("I31").Value = 3979.2727
$ws.Range("J31").Value = 5822.3335
$ws.Range("K31").Value = 3979.2727
$ws.Range("L31").Value = 5822.3335
$ws.Range("M31").Value = -3684.2727
$ws.Range("N31").Value = -6412.3335
$ws.Range("H34").Value = 4629.7646
$ws.Range("I34").Value = 3979.2727
$ws.Range("J34").Value = 5822.3335
$ws.Range("K34").Value = 3979.2727
$ws.Range("L34").Value = 5822.3335
$ws.Range("M34").Value = -3777.2727
$ws.Range("N34").Value = -6226.3335
$ws.Range("H86").Value = 12310.3125
$ws.Range("I86").Value = 5351.2
$ws.Range("K86").Value = 5351.2
$ws.Range("M86").Value = -4228.2
$ws.Range("H89").Value = 12310.3125
$ws.Range("I89").Value = 5351.2
$ws.Range("K89").Value = 26756
$ws.Range("M89").Value = -21140
$ws.Range("H99").Value = 6751333.5
$ws.Range("I99").Value = 10002000
$ws.Range("K99").Value = 10002000
$ws.Range("M99").Value = -10000502
$ws.Range("H105").Value = 2023.5555
$ws.Range("I105").Value = 1435.5
$ws.Range("J105").Value = 3199.6667
$ws.Range("K105").Value = 1435.5
$ws.Range("L105").Value = 3199.6667
$ws.Range("M105").Value = 311.5
$ws.Range("N105").Value = -6693.6667
$ws.Range("H126").Value = 6751333.5
$ws.Range("I126").Value = 10002000
$ws.Range("K126").Value = 30006000
$ws.Range("M126").Value = -30003530
$ws.Range("H132").Value = 9805830
$ws.Range("J132").Value = 17545606
$ws.Range("L132").Value = 52636818
$ws.Range("N132").Value = -52641878
$ws.Range("H134").Value = 2671.0188
$ws.Range("I134").Value = 1819.625
$ws.Range("K134").Value = 5458.875
$ws.Range("M134").Value = -2923.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2200
$ws.Range("J25").Value = 2925
$ws.Range("L25").Value = 8775
$ws.Range("N25").Value = -9113
$ws.Range("H30").Value = 2200
$ws.Range("J30").Value = 2925
$ws.Range("L30").Value = 8775
$ws.Range("N30").Value = -8979
$ws.Range("H44").Value = 111117336
$ws.Range("I44").Value = 261
$ws.Range("K44").Value = 783
$ws.Range("M44").Value = -385
$ws.Range("H87").Value = 15694.929
$ws.Range("J87").Value = 24890.285
$ws.Range("L87").Value = 74670.855
$ws.Range("N87").Value = -77166.855
$ws.Range("H90").Value = 15694.929
$ws.Range("J90").Value = 24890.285
$ws.Range("L90").Value = 224012.565
$ws.Range("N90").Value = -236492.565
$ws.Range("H131").Value = 25383522
$ws.Range("J131").Value = 28035272
$ws.Range("L131").Value = 84105816
$ws.Range("N131").Value = -84115896
$ws.Range("H137").Value = 9715.066000000001
$ws.Range("I137").Value = 7886
$ws.Range("J137").Value = 10629.6
$ws.Range("K137").Value = 23658
$ws.Range("L137").Value = 31888.8
$ws.Range("M137").Value = -18558
$ws.Range("N137").Value = -42088.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12721.366
$ws.Range("J70").Value = 12767.25
$ws.Range("L70").Value = 12767.25
$ws.Range("N70").Value = -13307.25
$ws.Range("H73").Value = 12721.366
$ws.Range("J73").Value = 12767.25
$ws.Range("L73").Value = 12767.25
$ws.Range("N73").Value = -14639.25
$ws.Range("H102").Value = 125001350
$ws.Range("I102").Value = 125001350
$ws.Range("K102").Value = 125001350
$ws.Range("M102").Value = -124999728

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3161.457
$ws.Range("I22").Value = 1377.1666
$ws.Range("J22").Value = 5050.706
$ws.Range("K22").Value = 1377.1666
$ws.Range("L22").Value = 5050.706
$ws.Range("M22").Value = -1082.1666
$ws.Range("N22").Value = -5640.706
$ws.Range("H27").Value = 3161.457
$ws.Range("I27").Value = 1377.1666
$ws.Range("J27").Value = 5050.706
$ws.Range("K27").Value = 1377.1666
$ws.Range("L27").Value = 5050.706
$ws.Range("M27").Value = -1270.1666
$ws.Range("N27").Value = -5264.706
$ws.Range("H40").Value = 103066
$ws.Range("I40").Value = 103066
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 103066
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -102930
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 8026.4287
$ws.Range("I132").Value = 4338.05
$ws.Range("K132").Value = 13014.15
$ws.Range("M132").Value = -10484.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6810.643
$ws.Range("J122").Value = 1977
$ws.Range("L122").Value = 5931
$ws.Range("N122").Value = -10831
